# Updated symbol list on Thu Feb  9 23:46:31 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with the latest scraped quotes.
# Values are plain text in the source workbook (t="inlineStr"), e.g. "0.1010"
# or "-6.56%", so force the Text number format before writing each cell -
# otherwise Excel auto-converts numeric-looking text into a real number/
# percentage and silently drops meaningful trailing zeros (e.g. "0.1010" -> 0.101).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.88"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-6.56%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-12.38%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.992"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.26%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07724"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-7.79%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.264"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.59%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.599"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-18.43%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9143"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.06%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1010"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.47%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1726"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-10.10%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09006"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.79%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04442"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.59%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.056"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-15.67%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1058"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.38%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001274"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.46%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005648"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.71%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.360"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.18%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.591"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.32%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3367"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.35%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.85%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2862"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.31%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04139"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.89%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001199"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.15%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004078"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-8.51%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001224"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.85%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002989"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.30%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02344"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-13.67%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05123"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-9.15%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007982"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.55%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1327"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.03%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007317"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.12%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002014"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.69%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008021"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.33%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3307"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.76%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006696"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.76%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.30%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003413"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.20%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004115"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "16.49%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002108"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.30%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.30%"
